$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Active cases"

$data = @(
    @("3398 BlueCross Elly Kay Mordialloc", 29),
    @("4257 BlueCross The Gables Camberwell", 16),
    @("44404 Castlemaine North Primary School Castlemaine", 10),
    @("4479 Whittlesea Lodge Whittlesea", 15),
    @("45168 Ranfurly Primary School Mildura", 11),
    @("52777 Mirripoa Primary School Mount DuneedSchool Camp", 14),
    @("Alfred Health The Alfred Hospital Melbourne", 14),
    @("Billboards The Venue Melbourne", 10),
    @("Confirmed Omicron Sircuit Bar Fitzroy", 12),
    @("Confirmed Omicron Variant The Peel Hotel Collingwood", 11),
    @("Feathertop Chalet Harrietville", 17),
    @("Melbourne Cricket Ground (MCG)", 12),
    @("Monash Health Dandenong Hospital Dandenong", 12),
    @("Ms Collins Melbourne", 40),
    @("St Pauls Cathedral", 60),
    @("St Vincents Hospital Melbourne Emergency Department Fitzroy", 30),
    @("St. Vincent's Hospital Melbourne Fitzroy", 24),
    @("The Emerson Rooftop Bar and Club South Yarra", 14),
    @("The Hatter and the Hare Bayswater", 22),
    @("The Royal Children's Hospital Melbourne Emergency Department Parkville", 13),
    @("Western Health Sunshine Hospital Emergency Department St Albans", 12)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
